$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of J2:J11
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14-17: summary labels and stats
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style the B14:B17 cells - bold, size 12, vertical center.
# Build the format once on a scratch cell, then paste-special (formats only)
# onto the target range so we don't fragment the style table with one new
# cellXfs entry per property assignment.
$scratch = $ws.Range("AA1")
$scratch.Font.Bold = $true
$scratch.Font.Size = 12
$scratch.VerticalAlignment = -4108
$scratch.Copy()
$ws.Range("B14:B17").PasteSpecial(-4122)
$scratch.Clear()

[void]$ws.Range("A14:B17").Select()

# Page setup (paper size / orientation) matching the re-saved workbook.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
